# feat: add 2022-Q4 data
#
# 1) Insert a new worksheet "2022-Q4" right after "总计", populated with
#    the new quarter's fund-holding rows (this pushes 2022-Q1 / 2021-Q4 /
#    2020-Q4 one slot later in tab order, matching the diff's sheetId/rId
#    renumbering).
# 2) Prepend a corresponding "2022-Q4" row to the "总计" summary sheet,
#    shifting the existing 2022-Q1 / 2021-Q4 / 2020-Q4 rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "总计" summary sheet: shift rows 2-4 down to 3-5, write new row 2
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.78

$total.Range("B3").Value = "2022-Q1"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.28

$total.Range("B4").Value = "2021-Q4"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.05

# Row 5 is brand new - clone the formatting of column A from row 4 first
# (same-sheet copy carries the cell style, s="2", over) before writing
# its value.
$total.Range("A4").Copy()
$total.Range("A5").PasteSpecial(-4122)
$total.Range("A5").Value = 3
$total.Range("B5").Value = "2020-Q4"
$total.Range("C5").Value = 1
$total.Range("D5").Value = 0.02

# ---------------------------------------------------------------------
# 2. New "2022-Q4" sheet: duplicate "2022-Q1" (an existing quarter sheet)
#    so the header/column-A cell styles come along for free, then
#    rename it and overwrite the cell values with the new quarter's data.
# ---------------------------------------------------------------------
$src = $wb.Worksheets.Item("2022-Q1")
$src.Copy($null, $total)
$newSheet = $wb.Worksheets.Item("2022-Q1 (2)")
$newSheet.Name = "2022-Q4"

# The source sheet only has 2 data rows; we need 3, so add row 4 by
# cloning row 2's column-A format (same-sheet copy keeps the s="2" style).
$newSheet.Range("A2").Copy()
$newSheet.Range("A4").PasteSpecial(-4122)

# Row 2
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'009805"
$newSheet.Range("C2").Value = "国泰医药健康股票A"
$newSheet.Range("D2").Value = "'10.35"
$newSheet.Range("E2").Value = "'93.35"
$newSheet.Range("F2").Value = "'6.89"
$newSheet.Range("G2").Value = "'0.7131"
$newSheet.Range("H2").Value = 7

# Row 3
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'011326"
$newSheet.Range("C3").Value = "国泰医药健康股票C"
$newSheet.Range("D3").Value = "'0.79"
$newSheet.Range("E3").Value = "'93.35"
$newSheet.Range("F3").Value = "'6.89"
$newSheet.Range("G3").Value = "'0.0544"
$newSheet.Range("H3").Value = 7

# Row 4
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'530016"
$newSheet.Range("C4").Value = "建信恒稳价值混合"
$newSheet.Range("D4").Value = "'0.54"
$newSheet.Range("E4").Value = "'52.32"
$newSheet.Range("F4").Value = "'3.09"
$newSheet.Range("G4").Value = "'0.0167"
$newSheet.Range("H4").Value = 6

# The leading apostrophes above force the numeric-looking strings to stay
# text, but they also stamp the cells with a quote-prefix style. Re-paste
# formatting from an always-plain text cell (C2, a fund-name string that
# was never apostrophe-prefixed) over those columns to strip that back off
# and match the unstyled text cells used elsewhere in the workbook.
$newSheet.Range("C2").Copy()
$newSheet.Range("B2:B4").PasteSpecial(-4122)
$newSheet.Range("D2:G4").PasteSpecial(-4122)

# Restore the originally active tab (last sheet, "2020-Q4") since creating/
# copying sheets above shifted the active selection to the new sheet.
$wb.Worksheets.Item("2020-Q4").Activate()
